$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns are treated as text so values like
# "1.00" or "65.18" are not auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.044.73"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.653.27"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").Value = "0.0614"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "20.06"
$ws.Range("E10").Value = "  +3.78%  "
$ws.Range("D11").Value = "0.0875"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "1.886.67"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "1.659.23"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "65.18"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "27.028.50"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "236.16"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "7.75"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "0.0₃0729"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "4.42"
$ws.Range("E22").Value = "  +3.63%  "
$ws.Range("E23").Value = "  +3.61%  "
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "145.35"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "7.08"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("D33").Value = "1.528.84"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +7.06%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "0.576"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.889"
$ws.Range("E38").Value = "  +7.86%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "5.92"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("D43").Value = "65.17"
$ws.Range("E43").Value = "  +6.85%  "
$ws.Range("D44").Value = "1.793.84"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("D46").Value = "0.914"
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "89.78"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").Value = "1.51"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "0.0506"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "0.0974"
$ws.Range("E51").Value = "  +0.92%  "
